$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("C2").Value = "chromegrid"
